$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(110, 3).Value = 16869
$ws.Cells.Item(110, 5).Value = 25930135

$ws.Cells.Item(115, 3).Value = 17549
$ws.Cells.Item(115, 5).Value = 38605108

$ws.Cells.Item(117, 3).Value = 19702
$ws.Cells.Item(117, 5).Value = 56430548

$ws.Cells.Item(121, 3).Value = 5963
$ws.Cells.Item(121, 5).Value = 11516074

$ws.Cells.Item(122, 3).Value = 9693
$ws.Cells.Item(122, 5).Value = 31946051

$ws.Cells.Item(134, 3).Value = 5673
$ws.Cells.Item(134, 5).Value = 17148182

$ws.Cells.Item(138, 3).Value = 2838
$ws.Cells.Item(138, 5).Value = 6585048

$ws.Cells.Item(139, 3).Value = 3317
$ws.Cells.Item(139, 5).Value = 9227172

$ws.Cells.Item(164, 3).Value = 50565
$ws.Cells.Item(164, 5).Value = 168372196

$ws.Cells.Item(168, 3).Value = 284920
$ws.Cells.Item(168, 5).Value = 1208341610

$ws.Cells.Item(169, 3).Value = 562561
$ws.Cells.Item(169, 5).Value = 1284323670

$ws.Cells.Item(170, 3).Value = 367272
$ws.Cells.Item(170, 5).Value = 2844005381

$ws.Cells.Item(171, 3).Value = 115103
$ws.Cells.Item(171, 5).Value = 444720413

$ws.Cells.Item(173, 3).Value = 54383
$ws.Cells.Item(173, 5).Value = 151846897

$ws.Cells.Item(174, 3).Value = 357160
$ws.Cells.Item(174, 5).Value = 1016615252

$ws.Cells.Item(175, 3).Value = 125508
$ws.Cells.Item(175, 5).Value = 811570009

$ws.Cells.Item(177, 3).Value = 96746
$ws.Cells.Item(177, 5).Value = 174705213

$ws.Cells.Item(179, 3).Value = 235655
$ws.Cells.Item(179, 5).Value = 812084503

$ws.Cells.Item(186, 3).Value = 21933
$ws.Cells.Item(186, 5).Value = 40056540

$ws.Cells.Item(188, 3).Value = 19703
$ws.Cells.Item(188, 5).Value = 66031872

$ws.Cells.Item(196, 3).Value = 7400
$ws.Cells.Item(196, 5).Value = 20652188

$ws.Cells.Item(198, 3).Value = 4508
$ws.Cells.Item(198, 5).Value = 5999357

$ws.Cells.Item(199, 3).Value = 4156
$ws.Cells.Item(199, 5).Value = 9036256

$ws.Cells.Item(203, 3).Value = 13101
$ws.Cells.Item(203, 5).Value = 32996553

$ws.Cells.Item(204, 3).Value = 4754
$ws.Cells.Item(204, 5).Value = 11644170

$ws.Cells.Item(205, 3).Value = 11123
$ws.Cells.Item(205, 5).Value = 44077924

$ws.Cells.Item(209, 3).Value = 5363
$ws.Cells.Item(209, 5).Value = 12210515

$ws.Cells.Item(211, 3).Value = 2863
$ws.Cells.Item(211, 5).Value = 4380689

$ws.Cells.Item(213, 3).Value = 3633
$ws.Cells.Item(213, 5).Value = 11097853

$ws.Cells.Item(214, 3).Value = 6172
$ws.Cells.Item(214, 5).Value = 11075372

$ws.Cells.Item(241, 3).Value = 2583
$ws.Cells.Item(241, 5).Value = 7741099

$ws.Cells.Item(267, 3).Value = 84974
$ws.Cells.Item(267, 5).Value = 156518758

$ws.Cells.Item(295, 3).Value = 91332
$ws.Cells.Item(295, 5).Value = 552911368

$ws.Cells.Item(317, 3).Value = 103579
$ws.Cells.Item(317, 5).Value = 303078087

$ws.Cells.Item(320, 3).Value = 67241
$ws.Cells.Item(320, 5).Value = 124554315

$ws.Cells.Item(322, 3).Value = 81161
$ws.Cells.Item(322, 4).Value = 9703
$ws.Cells.Item(322, 5).Value = 254527631
